# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G (header "K") values are recalculated/rewritten with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-19 (row index -> new value)
$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    6  = 3
    7  = 2
    8  = 1
    9  = 2
    10 = 0
    11 = 2
    12 = 1
    13 = 3
    14 = 1
    15 = 3
    16 = 2
    17 = 3
    18 = 3
    19 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
